$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("J6").Value = 2.8
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 3.5
$ws.Range("N6").Value = 8
$ws.Range("P6").Value = 2.92
$ws.Range("W6").Value = 7.7
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 9
$ws.Range("AA6").Value = 18.5
$ws.Range("AB6").Value = 28
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 16
$ws.Range("AO6").Value = 11.5
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 45
$ws.Range("AR6").Value = 70
$ws.Range("AT6").Value = 2.62
$ws.Range("AX6").Value = 16
$ws.Range("AY6").Value = 21
$ws.Range("BA6").Value = 100
$ws.Range("G12").Value = 2.3
$ws.Range("I12").Value = 2.9
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
$ws.Range("Y12").Value = 9.5
$ws.Range("AC12").Value = 8.5
$ws.Range("AD12").Value = 6.5
$ws.Range("AW12").Value = 5
$ws.Range("AY12").Value = 29
$ws.Range("G17").Value = 2.05
$ws.Range("I17").Value = 3.4
$ws.Range("J17").Value = 2.63
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("Z17").Value = 19
$ws.Range("AH17").Value = 19
$ws.Range("AI17").Value = 12
$ws.Range("AO17").Value = 11
$ws.Range("AX17").Value = 17
$ws.Range("G25").Value = 2.1
$ws.Range("I25").Value = 3
$ws.Range("J25").Value = 2.7
$ws.Range("K25").Value = 2.2
$ws.Range("L25").Value = 3.55
$ws.Range("M25").Value = 1.05
$ws.Range("N25").Value = 8
$ws.Range("P25").Value = 3.6
$ws.Range("Q25").Value = 1.75
$ws.Range("R25").Value = 2.02
$ws.Range("S25").Value = 1.36
$ws.Range("T25").Value = 2.92
$ws.Range("V25").Value = 2.15
$ws.Range("W25").Value = 8.75
$ws.Range("X25").Value = 11.25
$ws.Range("Z25").Value = 20
$ws.Range("AA25").Value = 16
$ws.Range("AC25").Value = 8
$ws.Range("AE25").Value = 13
$ws.Range("AG25").Value = 10.75
$ws.Range("AH25").Value = 17
$ws.Range("AI25").Value = 10.75
$ws.Range("AJ25").Value = 37
$ws.Range("AK25").Value = 24
$ws.Range("AL25").Value = 29
$ws.Range("AM25").Value = 350
$ws.Range("AN25").Value = 4.15
$ws.Range("AO25").Value = 10.75
$ws.Range("AT25").Value = 2.92
$ws.Range("AU25").Value = 6.9
$ws.Range("AW25").Value = 5.1
$ws.Range("AX25").Value = 16
$ws.Range("AY25").Value = 22
$ws.Range("AZ25").Value = 75
$ws.Range("BA25").Value = 100
$ws.Range("K26").Value = 1.85
$ws.Range("N26").Value = 5.7
$ws.Range("O26").Value = 1.5
$ws.Range("P26").Value = 2.25
$ws.Range("Q26").Value = 2.45
$ws.Range("R26").Value = 1.42
$ws.Range("S26").Value = 1.53
$ws.Range("T26").Value = 2.18
$ws.Range("U26").Value = 1.98
$ws.Range("V26").Value = 1.65
$ws.Range("W26").Value = 6.7
$ws.Range("Y26").Value = 11
$ws.Range("AC26").Value = 5.9
$ws.Range("AD26").Value = 5.4
$ws.Range("AE26").Value = 17
$ws.Range("AG26").Value = 6.4
$ws.Range("AH26").Value = 12
$ws.Range("AI26").Value = 10.25
$ws.Range("AN26").Value = 4.6
$ws.Range("AP26").Value = 28
$ws.Range("AT26").Value = 2.15
$ws.Range("AU26").Value = 7.4
$ws.Range("AV26").Value = 80
$ws.Range("AY26").Value = 26
$ws.Range("BA26").Value = 120
$ws.Range("BB26").Value = 400
